# TABLA DE CONTROL DE ACTIVIDADES - add SEMANA 27/28/29 entries (sección 5 GIF work)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOVIEMBRE 2024")

# ---------------------------------------------------------------------------
# 1) Remove the stray trailing date row (old B247=45837 with empty activity),
#    which collapses the old SUM row (248 -> 247) and auto-fixes its range.
# ---------------------------------------------------------------------------
$ws.Rows.Item(247).Delete()

# ---------------------------------------------------------------------------
# 2) The week block starting at row 241 is now "SEMANA 26" (was "SEMANA 25").
# ---------------------------------------------------------------------------
$ws.Range("B241").Value = "SEMANA 26"

# ---------------------------------------------------------------------------
# Helper template for a week block: header row (merged B:F), column-titles
# row, one-or-more data rows, then a SUM row in column G. We clone
# formatting from the existing "SEMANA 22" block (rows 213-216), which has
# exactly this one-data-row shape, then fill in the new content.
# ---------------------------------------------------------------------------

# --- SEMANA 27 : rows 249-252 ---------------------------------------------
$ws.Range("B213:F215").Copy()
$ws.Range("B249").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("G238").Copy()
$ws.Range("G252").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B249").Value = "SEMANA 27"
$ws.Range("B249:F249").Merge()

$ws.Range("B250").Value = "FECHA"
$ws.Range("C250").Value = "ACTIVIDAD"
$ws.Range("D250").Value = "OBSERVACIONES"
$ws.Range("E250").Value = "HORARIO"
$ws.Range("F250").Value = "HORAS CUBIERTAS"

$ws.Range("B251").Value = 45841
$ws.Range("C251").Value = "Se realizarón los cambion responsivos para la sección 5"
$ws.Range("D251").Value = "Se empezarón los cambios responsivos parara la sección 5"
$ws.Range("E251").Value = "11:00-14:00"
$ws.Range("F251").Value = 3
$ws.Rows.Item(251).RowHeight = 45

$ws.Range("G252").Formula = "=SUM(F251:F251)"

# --- SEMANA 28 : rows 255-258 ---------------------------------------------
$ws.Range("B213:F215").Copy()
$ws.Range("B255").PasteSpecial(-4122)
$ws.Range("G238").Copy()
$ws.Range("G258").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B255").Value = "SEMANA 28"
$ws.Range("B255:F255").Merge()

$ws.Range("B256").Value = "FECHA"
$ws.Range("C256").Value = "ACTIVIDAD"
$ws.Range("D256").Value = "OBSERVACIONES"
$ws.Range("E256").Value = "HORARIO"
$ws.Range("F256").Value = "HORAS CUBIERTAS"

$ws.Range("B257").Value = 45862
$ws.Range("C257").Value = "Se terminó la funcionalidad para la sección 5, se corrigierón detalles y se empezo a trabajar en la sección 6"
$ws.Range("D257").Value = "Se terminó la funcionalidad para la sección 5, se corrigieron detalles y se empezo la sección 6, para la planta de tratamiento AR. Se presentarón problemas con los GIFs"
$ws.Range("E257").Value = "7:00-11:30, 15:30-20:00"
$ws.Range("F257").Value = 9
$ws.Rows.Item(257).RowHeight = 90

$ws.Range("G258").Formula = "=SUM(F257:F257)"

# --- SEMANA 29 : rows 261-264 ---------------------------------------------
$ws.Range("B213:F215").Copy()
$ws.Range("B261").PasteSpecial(-4122)
$ws.Range("G238").Copy()
$ws.Range("G264").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B261").Value = "SEMANA 29"
$ws.Range("B261:F261").Merge()

$ws.Range("B262").Value = "FECHA"
$ws.Range("C262").Value = "ACTIVIDAD"
$ws.Range("D262").Value = "OBSERVACIONES"
$ws.Range("E262").Value = "HORARIO"
$ws.Range("F262").Value = "HORAS CUBIERTAS"

$ws.Range("B263").Value = 45866
$ws.Range("C263").Value = "Se creó una función para el control del los GIFs, considerando su reinicio y optimización"
$ws.Range("D263").Value = "Se terminó una función para el control de GIFs de la sección 5, para la planta de tratamiento AR"
$ws.Range("E263").Value = "8:30-13:00, 15:30-20:00"
$ws.Range("F263").Value = 9
$ws.Rows.Item(263).RowHeight = 60

$ws.Range("G264").Formula = "=SUM(F263:F263)"

# ---------------------------------------------------------------------------
# 3) Fix up the total-hours summary formula: the author trimmed the big
#    SUM range from G12:G1200 down to G12:G1193.
# ---------------------------------------------------------------------------
$ws.Range("J6").Formula = "=SUM(G12:G1193)"

# ---------------------------------------------------------------------------
# 4) Update the visible selection to mirror the author's final cursor spot.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("D271").Select()
